# Updates cryptos.xlsx price/volume figures and reorders a couple of coin rows,
# matching the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.269.85"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.895.44"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'246.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'0.691"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.30%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'40.40"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.21%  "
$ws.Range("D9").Value = "'0.347"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").Value = "'51.93"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.48%  "
$ws.Range("D11").Value = "'0.0721"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "2.172.05"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "'12.51"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "'0.707"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "1.888.44"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "'4.82"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "35.260.06"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "'72.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "0.0₃0818"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'240.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").Value = "'12.73"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'2.32"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.59%  "
$ws.Range("D27").Value = "'167.99"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("D28").Value = "'8.55"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").Value = "'19.02"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.42%  "
$ws.Range("D30").Value = "'0.130"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("D32").Value = "'4.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").Value = "'0.0569"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +6.88%  "
$ws.Range("D36").Value = "'4.12"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").Value = "'0.910"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.38%  "
$ws.Range("D38").Value = "'1.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.01%  "
$ws.Range("D39").Value = "'2.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'0.0656"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.78%  "
$ws.Range("D41").Value = "'94.64"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'16.44"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.56%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.09"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").Value = "'0.0207"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "1.351.28"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").Value = "'2.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "Gas"
$ws.Range("C48").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D48").Value = "'12.51"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.78"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "'45.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.41%  "
$ws.Range("D51").Value = "'6.48"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.09%  "
